# "subindo atualizacao do timescore e execução de experiencia"
#
# Updates the TimeScore example workbook:
#   - changes the sample input years on row 5 of the TimeScore sheet
#   - adds a new "doubled time score" column (P) for row 5
#   - adds a new experiment block on row 8 (average hm / alternate score)
#   - leaves row 6 as a blank spacer row
#   - makes the TimeScore sheet the active/selected sheet with D5 selected
#   - leaves "Ideias do Professor" with its D19 selection, no longer active

$wb = $excel.ActiveWorkbook
$wsTime = $wb.Worksheets.Item("TimeScore")
$wsProf = $wb.Worksheets.Item("Ideias do Professor")

# --- Row 5: update the sample input values (years) ---
$wsTime.Range("D5").Value = 2004
$wsTime.Range("E5").Value = 2004
$wsTime.Range("F5").Value = 2004
$wsTime.Range("H5").Value = 2005
# (G5 stays 0.5, formulas I5/J5/K5/L5/M5/N5 are untouched and simply
#  recalculate from the new inputs)

# --- New column P on row 5: doubled time score ---
$wsTime.Range("P5").Formula = "=N5+N5"

# --- Row 6 is left completely empty, just used as a spacer ---

# --- Row 8: new "alternate" experiment using the average of B5/C5 ---
$wsTime.Range("I8").Formula = "=(B5+C5)/2"
$wsTime.Range("L8").Formula = "=I8*K5"
$wsTime.Range("M8").Value = 2
$wsTime.Range("N8").Formula = "=L8/M8"

# --- Match the new row heights used for the edited/added rows ---
$wsTime.Rows.Item(5).RowHeight = 12.8
$wsTime.Rows.Item(6).RowHeight = 12.8
$wsTime.Rows.Item(8).RowHeight = 12.8

# --- View/selection state ---
# Select D19 on the professor sheet first (keeps that selection stored
# for that sheet) then make TimeScore the active sheet with D5 selected,
# matching the tabSelected/activeTab flip in the diff.
$wsProf.Range("D19").Select() | Out-Null

$wsTime.Activate() | Out-Null
$wsTime.Range("D5").Select() | Out-Null
